# Applies a reordering of several observation records within the
# "Artfynd" sheet. The records (full rows of data) are cyclically
# rotated among specific row slots, while shared/contextual columns
# (location, date, observer, etc.) stay untouched because they are
# identical across all affected rows.
#
# Only the columns that actually carry record-specific data are
# touched: A, B, D, E, F, G, H, K, L, M, N, Q, R, Z, AB, AC.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","K","L","M","N","Q","R","Z","AB","AC")

# Explicit destination-row -> source-row map: the destination row
# (key) receives the record data that originally lived in the source
# row (value), matching the diff exactly:
#   row 4  <- old row 5     row 5  <- old row 6     row 6  <- old row 4
#   row 7  <- old row 8     row 8  <- old row 9     row 9  <- old row 10   row 10 <- old row 7
#   row 12 <- old row 13    row 13 <- old row 12
#   row 29 <- old row 30    row 30 <- old row 31     row 31 <- old row 29
#   row 39 <- old row 40    row 40 <- old row 39
$destToSrc = @{
    4  = 5
    5  = 6
    6  = 4
    7  = 8
    8  = 9
    9  = 10
    10 = 7
    12 = 13
    13 = 12
    29 = 30
    30 = 31
    31 = 29
    39 = 40
    40 = 39
}

# Snapshot all relevant cell values first, before making any changes,
# so that later writes never clobber data still needed as a source.
$snapshot = @{}
foreach ($r in $destToSrc.Keys) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

foreach ($destRow in $destToSrc.Keys) {
    $srcRow = $destToSrc[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcData[$c]
    }
}
